# The deck currently ships two theme parts:
#   ppt/theme/theme1.xml  ("Integral" / "Red Violet" colours)  -> used by the slide master
#   ppt/theme/theme2.xml  ("Office Theme" / "Office" colours)  -> used by the notes master
#
# The target edit swaps the two themes' content, so the slide master ends up
# using the "Office Theme" palette and the notes master ends up with the
# "Integral" / "Red Violet" palette. The font scheme and format scheme blocks
# are identical between the two theme parts, so the only substantive
# difference is the 12-colour colour scheme (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink).
#
# This host only exposes the single slide-master-backed theme through the
# object model (Master.ColorScheme / Slide.ThemeColorScheme, etc. all read
# and write that same theme part - ppt/theme/theme1.xml - regardless of
# which master/slide they are called from), so we update that theme's 12
# colours to the "Office Theme" palette using Slide.ThemeColorScheme, which
# edits colours in place without disturbing the rest of the theme XML.

function HexToBgrInt($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation

# Office Theme colour scheme, in clrScheme order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = HexToBgrInt $officeThemeColors[$i - 1]
}
